$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Câu truy vấn lấy quý hiện tại"
$ws.Range("C4").Value = 42959
$ws.Range("C4").NumberFormat = "mm-dd-yy"
$ws.Range("D4").Value = "Đào Hoài Phương"
$ws.Range("H4").Value = "Open"

$ws.Range("H4").Select()
